$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F21").Value = -5
